$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 45678
$ws.Range("J108").Value = 45678
$ws.Range("L108").Value = 45678
$ws.Range("N108").Value = -53358
$ws.Range("H114").Value = 46012.4
$ws.Range("J114").Value = 46012.4
$ws.Range("L114").Value = 46012.4
$ws.Range("N114").Value = -54690.4
$ws.Range("H120").Value = 49722
$ws.Range("J120").Value = 49722
$ws.Range("L120").Value = 49722
$ws.Range("N120").Value = -59398

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1128.9457
$ws.Range("I61").Value = 1002.3125
$ws.Range("J61").Value = 1418.3928
$ws.Range("K61").Value = 1002.3125
$ws.Range("L61").Value = 1418.3928
$ws.Range("M61").Value = -790.3125
$ws.Range("N61").Value = -1842.3928
$ws.Range("H74").Value = 1954.8551
$ws.Range("I74").Value = 1964.8077
$ws.Range("J74").Value = 1924.4117
$ws.Range("K74").Value = 1964.8077
$ws.Range("L74").Value = 1924.4117
$ws.Range("M74").Value = -1090.8077
$ws.Range("N74").Value = -3672.4117
$ws.Range("H77").Value = 1954.8551
$ws.Range("I77").Value = 1964.8077
$ws.Range("J77").Value = 1924.4117
$ws.Range("K77").Value = 9824.038500000001
$ws.Range("L77").Value = 9622.058500000001
$ws.Range("M77").Value = -5456.038500000001
$ws.Range("N77").Value = -18358.0585
$ws.Range("H102").Value = 12062
$ws.Range("I102").Value = 1219.909
$ws.Range("K102").Value = 1219.909
$ws.Range("M102").Value = 402.0909999999999
$ws.Range("H111").Value = 46537.332
$ws.Range("J111").Value = 46537.332
$ws.Range("L111").Value = 46537.332
$ws.Range("N111").Value = -54717.332
$ws.Range("H121").Value = 44987
$ws.Range("J121").Value = 44987
$ws.Range("L121").Value = 44987
$ws.Range("N121").Value = -48481
$ws.Range("H136").Value = 1128.9457
$ws.Range("I136").Value = 1002.3125
$ws.Range("J136").Value = 1418.3928
$ws.Range("K136").Value = 3006.9375
$ws.Range("L136").Value = 4255.178400000001
$ws.Range("M136").Value = -456.9375
$ws.Range("N136").Value = -9355.178400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 47684
$ws.Range("J108").Value = 47684
$ws.Range("L108").Value = 47684
$ws.Range("N108").Value = -55364
$ws.Range("H112").Value = 47469
$ws.Range("J112").Value = 47469
$ws.Range("L112").Value = 47469
$ws.Range("N112").Value = -50423
$ws.Range("H116").Value = 45734
$ws.Range("J116").Value = 45734
$ws.Range("L116").Value = 45734
$ws.Range("N116").Value = -54912

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H110").Value = 37042
$ws.Range("J110").Value = 37042
$ws.Range("L110").Value = 37042
$ws.Range("N110").Value = -45222
$ws.Range("H111").Value = 48702
$ws.Range("J111").Value = 48702
$ws.Range("L111").Value = 48702
$ws.Range("N111").Value = -56882
$ws.Range("H116").Value = 49364.5
$ws.Range("J116").Value = 49364.5
$ws.Range("L116").Value = 49364.5
$ws.Range("N116").Value = -58542.5
$ws.Range("H118").Value = 48074
$ws.Range("J118").Value = 48074
$ws.Range("L118").Value = 48074
$ws.Range("N118").Value = -51388
$ws.Range("H119").Value = 46504.332
$ws.Range("J119").Value = 46504.332
$ws.Range("L119").Value = 46504.332
$ws.Range("N119").Value = -56180.332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 41544.6
$ws.Range("J114").Value = 41544.6
$ws.Range("L114").Value = 41544.6
$ws.Range("N114").Value = -50222.6
$ws.Range("H116").Value = 38998
$ws.Range("J116").Value = 38998
$ws.Range("L116").Value = 38998
$ws.Range("N116").Value = -48176
$ws.Range("H119").Value = 48753
$ws.Range("J119").Value = 48753
$ws.Range("L119").Value = 48753
$ws.Range("N119").Value = -58429
$ws.Range("H126").Value = 8523.5
$ws.Range("I126").Value = 15215.75
$ws.Range("K126").Value = 45647.25
$ws.Range("M126").Value = -43177.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 48622
$ws.Range("J108").Value = 48622
$ws.Range("L108").Value = 48622
$ws.Range("N108").Value = -56302
$ws.Range("H111").Value = 43744
$ws.Range("J111").Value = 43744
$ws.Range("L111").Value = 43744
$ws.Range("N111").Value = -51924
$ws.Range("H112").Value = 41970.668
$ws.Range("J112").Value = 41970.668
$ws.Range("L112").Value = 41970.668
$ws.Range("N112").Value = -44924.668
$ws.Range("H114").Value = 27593.6
$ws.Range("J114").Value = 27593.6
$ws.Range("L114").Value = 27593.6
$ws.Range("N114").Value = -36271.6
$ws.Range("H116").Value = 48586.668
$ws.Range("J116").Value = 48586.668
$ws.Range("L116").Value = 48586.668
$ws.Range("N116").Value = -57764.668
$ws.Range("H117").Value = 43380
$ws.Range("J117").Value = 43380
$ws.Range("L117").Value = 43380
$ws.Range("N117").Value = -52558
$ws.Range("H118").Value = 40350.25
$ws.Range("J118").Value = 40350.25
$ws.Range("L118").Value = 40350.25
$ws.Range("N118").Value = -43664.25
$ws.Range("H119").Value = 43194.668
$ws.Range("J119").Value = 43194.668
$ws.Range("L119").Value = 43194.668
$ws.Range("N119").Value = -52870.668
$ws.Range("H120").Value = 45835.6
$ws.Range("J120").Value = 45835.6
$ws.Range("L120").Value = 45835.6
$ws.Range("N120").Value = -55511.6
$ws.Range("H121").Value = 40084
$ws.Range("J121").Value = 40084
$ws.Range("L121").Value = 40084
$ws.Range("N121").Value = -43578

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 44626
$ws.Range("J108").Value = 44626
$ws.Range("L108").Value = 44626
$ws.Range("N108").Value = -52306
$ws.Range("H109").Value = 39377
$ws.Range("J109").Value = 39377
$ws.Range("L109").Value = 39377
$ws.Range("N109").Value = -42151
$ws.Range("H110").Value = 48636
$ws.Range("J110").Value = 48636
$ws.Range("L110").Value = 48636
$ws.Range("N110").Value = -56816
$ws.Range("H115").Value = 37932
$ws.Range("J115").Value = 37932
$ws.Range("L115").Value = 37932
$ws.Range("N115").Value = -41066
$ws.Range("H116").Value = 50000
$ws.Range("J116").Value = 50000
$ws.Range("L116").Value = 50000
$ws.Range("N116").Value = -59178
$ws.Range("H117").Value = 44075.5
$ws.Range("J117").Value = 44075.5
$ws.Range("L117").Value = 44075.5
$ws.Range("N117").Value = -53253.5
$ws.Range("H118").Value = 43384
$ws.Range("J118").Value = 43384
$ws.Range("L118").Value = 43384
$ws.Range("N118").Value = -46698
$ws.Range("H119").Value = 48698
$ws.Range("J119").Value = 48698
$ws.Range("L119").Value = 48698
$ws.Range("N119").Value = -58374
$ws.Range("H121").Value = 41916
$ws.Range("J121").Value = 41916
$ws.Range("L121").Value = 41916
$ws.Range("N121").Value = -45410
$ws.Range("H126").Value = 1963114.2
$ws.Range("J126").Value = 2454.9
$ws.Range("L126").Value = 7364.700000000001
$ws.Range("N126").Value = -12304.7
